# Apply the latest crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.555.08"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.519.77"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.26"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.79"
$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.515.91"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E10").Value = "  -3.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  +7.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.31"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.090.78"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.29"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "607.78"
$ws.Range("E17").Value = "  -2.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.518.66"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.631.14"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.47"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("E23").Value = "  -9.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.53"
$ws.Range("E24").Value = "  +2.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.63"
$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.73"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.57"
$ws.Range("E28").Value = "  -1.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.91"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.03"
$ws.Range("E30").Value = "  -2.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("E32").Value = "  -4.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "637.11"
$ws.Range("E33").Value = "  +11.98%  "

$ws.Range("E34").Value = "  -4.58%  "

$ws.Range("E35").Value = "  -2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.58"
$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0994"
$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.76"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").Value = "  +4.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.81"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.371.74"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0741"
$ws.Range("E44").Value = "  +4.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.310"
$ws.Range("E45").Value = "  -5.40%  "

$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.13"
$ws.Range("E47").Value = "  -3.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.56"
$ws.Range("E48").Value = "  -3.50%  "

$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.48"
$ws.Range("E50").Value = "  -2.11%  "
